$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 2 new rows before row 16 (pushes existing rows 16.. down to 18..) ---
$ws.Rows("16:17").Insert()

# --- Fill new row 17 first (RS485 to TTL), then row 16 (USB to RS485) so that
#     shared-string allocation order matches the authored workbook ---
$ws.Range("A17").Value = "RS485 to TTL"
$ws.Range("B17").Value = "RS485 shifter for Redox"
$ws.Range("C17").Value = "https://www.amazon.com/gp/product/B01FDD52Y2/ref=ppx_yo_dt_b_asin_title_o07_s00?ie=UTF8&psc=1"
$ws.Range("D17").Value = 7.19
$ws.Range("E17").Value = 1
$ws.Range("F17").Formula = "=D17*E17"

$ws.Range("A16").Value = "USB to RS485"
$ws.Range("B16").Value = "Useful to read Redox data from computer"
$ws.Range("C16").Value = "https://www.amazon.com/gp/product/B00NKAJGZM/ref=ppx_yo_dt_b_asin_title_o09_s00?ie=UTF8&psc=1"
$ws.Range("D16").Value = 6.99
$ws.Range("E16").Value = 1
$ws.Range("F16").Formula = "=D16*E16"

# C16/C17 use the hyperlink style (same look as other link cells in column C)
$ws.Range("C16").Style = $ws.Range("C15").Style
$ws.Range("C17").Style = $ws.Range("C15").Style

# --- Rebuild hyperlinks: row-insert does not auto-shift stored hyperlink
#     ranges, so clear everything and re-add at the correct (shifted) cells,
#     plus the two brand-new ones for the RS485 rows. ---
$ws.Range("A1").Hyperlinks.Delete() | Out-Null

$ws.Hyperlinks.Add($ws.Range("C50"), "https://www.automationdirect.com/adc/Shopping/Catalog/Wiring_Solutions/Micro_(M12)_Receptacles/7231-13501-9710050") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C49"), "https://www.smarthome.com/elk-w040a-m1-cable-to-db9-serial-ribbon-cable.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://microflx.com/products/microlink-hm-hart-protocol-modem-modbus-accumulator-rs-232-interface?variant=33428219782") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.procomsol.com/online_store/r_loop_250_ohm_hart_loop_resistor") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "https://www.sparkfun.com/products/298") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C17"), "https://www.amazon.com/gp/product/B01FDD52Y2/ref=ppx_yo_dt_b_asin_title_o07_s00?ie=UTF8&psc=1") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C16"), "https://www.amazon.com/gp/product/B00NKAJGZM/ref=ppx_yo_dt_b_asin_title_o09_s00?ie=UTF8&psc=1") | Out-Null

# --- Column H width (auto-fit picks up the "Total" header + $ totals that now
#     sit in that column) ---
$ws.Columns("H").ColumnWidth = 9.6

# --- Selection moves to the newly-added F16 cell ---
$ws.Range("F16").Select() | Out-Null

# --- Window was maximized when the file was last saved ---
$excel.ActiveWindow.WindowState = -4137 | Out-Null
